# Update import product detail table format.xlsx
# (import maks 10000 row berikutnya)
#
# The first data row (row 3: 40 / DL00023-1LUBCK / DK) is a duplicate/stale
# entry that needs to be skipped, so we delete it entirely. This shifts all
# the rows below it up by one (row 4 -> row 3, row 5 -> row 4, ... row 63 -> row 62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire 3rd row (the first data row after the two header rows).
$ws.Rows.Item(3).Delete()

# Reflect the cursor position left behind by the edit.
$ws.Range("C10").Select()
